# feat: Añadir cambiar contraseña en vMtoUsuarios
#
# Adds a new row describing the "$_SESSION [codDepartamentoEnCursoRest]"
# variable (the department code picked from the new "Rest" window, used by
# "Mi Api") to the "uso de la sesion" reference table on Hoja1, and moves
# the selection to the cell the author ended up on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# New row 15.
# ----------------------------------------------------------------------
$ws.Rows.Item(15).RowHeight = 29.85

# B15: "$_SESSION [codDepartamentoEnCursoRest]" with the variable name in
# bold + accent colour, matching every other "$_SESSION [...]" entry above
# it (e.g. B11 "$_SESSION [codDepartamentoEnCurso]").
$prefix = '$_SESSION ['
$varName = 'codDepartamentoEnCursoRest'
$suffix = ']'
$ws.Range("B15").Value = $prefix + $varName + $suffix

$ws.Range("B15").HorizontalAlignment = -4131
$ws.Range("B15").VerticalAlignment = -4108
$ws.Range("B15").Borders.LineStyle = 1

$prefixLen = $prefix.Length
$varLen = $varName.Length

$charsPrefix = $ws.Range("B15").Characters(1, $prefixLen)
$charsPrefix.Font.Size = 11
$charsPrefix.Font.Name = "Calibri"
$charsPrefix.Font.Color = 0

$charsVar = $ws.Range("B15").Characters($prefixLen + 1, $varLen)
$charsVar.Font.Bold = $true
$charsVar.Font.Size = 11
$charsVar.Font.Name = "Calibri"
$charsVar.Font.Color = 2315831

$charsSuffix = $ws.Range("B15").Characters($prefixLen + $varLen + 1, $suffix.Length)
$charsSuffix.Font.Size = 11
$charsSuffix.Font.Name = "Calibri"
$charsSuffix.Font.Color = 0

# C15: plain description text, wrapped, with the same border used by the
# rest of column C.
$ws.Range("C15").Value = 'Texto que guarda el codigo del departamento seleccionado desde la ventana Rest para su uso en la Mi Api.'
$ws.Range("C15").HorizontalAlignment = -4131
$ws.Range("C15").VerticalAlignment = -4108
$ws.Range("C15").WrapText = $true
$ws.Range("C15").Borders.LineStyle = 1

# ----------------------------------------------------------------------
# C13 / C14 end up sharing the same wrapped-border style as the rest of
# column C once the new row is added.
# ----------------------------------------------------------------------
$ws.Range("C13:C14").WrapText = $true

# ----------------------------------------------------------------------
# The author's cursor ended up on E14 when they finished editing.
# ----------------------------------------------------------------------
$ws.Range("E14").Select() | Out-Null

Write-Output "done"
